$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from row 16 column A down to new rows 17-19 so they keep the bold/bordered/centered style used for column A throughout the table
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update rows 10-19 with the averaged-intensity data, inserting the new Gaussian-Quadrature/Spiral schemes ---
# Row 10: Gaussian-Quadrature
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.045273451079792
$ws.Cells.Item(10, 4).Value = 0.8527933969363093
$ws.Cells.Item(10, 5).Value = 1.019274105624675
$ws.Cells.Item(10, 6).Value = 1.045273451079792
$ws.Cells.Item(10, 7).Value = 0.920361292935602
$ws.Cells.Item(10, 8).Value = 1.065469375697556
$ws.Cells.Item(10, 9).Value = 1.02580319011694
$ws.Cells.Item(10, 10).Value = 0.8527933969363093
$ws.Cells.Item(10, 11).Value = 0.9360337512804924
$ws.Cells.Item(10, 12).Value = 0.9906536011801422
$ws.Cells.Item(10, 13).Value = 0.9881624687318125

# Row 11: Spiral-90deg-10rot-5space
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 1.004854987927426
$ws.Cells.Item(11, 4).Value = 0.9434705161733614
$ws.Cells.Item(11, 5).Value = 1.008726020025434
$ws.Cells.Item(11, 6).Value = 1.004854987927426
$ws.Cells.Item(11, 7).Value = 0.9641373337193397
$ws.Cells.Item(11, 8).Value = 1.033240208862283
$ws.Cells.Item(11, 9).Value = 1.008413879685662
$ws.Cells.Item(11, 10).Value = 0.9434705161733614
$ws.Cells.Item(11, 11).Value = 0.9760982680993975
$ws.Cells.Item(11, 12).Value = 0.9904766280134117
$ws.Cells.Item(11, 13).Value = 0.993807157732251

# Row 12: Spiral-90deg-15rot-5space
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 1.004632827689991
$ws.Cells.Item(12, 4).Value = 0.9438456484956003
$ws.Cells.Item(12, 5).Value = 1.008713064532918
$ws.Cells.Item(12, 6).Value = 1.004632827689991
$ws.Cells.Item(12, 7).Value = 0.9644175866334675
$ws.Cells.Item(12, 8).Value = 1.033125190062284
$ws.Cells.Item(12, 9).Value = 1.008349371256998
$ws.Cells.Item(12, 10).Value = 0.9438456484956003
$ws.Cells.Item(12, 11).Value = 0.9762793565142591
$ws.Cells.Item(12, 12).Value = 0.9904560921021253
$ws.Cells.Item(12, 13).Value = 0.9938472814452098

# Row 13: Spiral-90deg-10rot-3space
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 1.004786411471297
$ws.Cells.Item(13, 4).Value = 0.9434399772312331
$ws.Cells.Item(13, 5).Value = 1.008765871892324
$ws.Cells.Item(13, 6).Value = 1.004786411471297
$ws.Cells.Item(13, 7).Value = 0.9642525903886964
$ws.Cells.Item(13, 8).Value = 1.033275957260047
$ws.Cells.Item(13, 9).Value = 1.008419894708838
$ws.Cells.Item(13, 10).Value = 0.9434399772312331
$ws.Cells.Item(13, 11).Value = 0.9761029245617785
$ws.Cells.Item(13, 12).Value = 0.9904446680165379
$ws.Cells.Item(13, 13).Value = 0.9938234504920725

# Row 14: NoRotation-tilt60deg
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 1.055627999999998
$ws.Cells.Item(14, 4).Value = 0.7055639999999994
$ws.Cells.Item(14, 5).Value = 1.080023999999999
$ws.Cells.Item(14, 6).Value = 1.055627999999998
$ws.Cells.Item(14, 7).Value = 0.8028560000000003
$ws.Cells.Item(14, 8).Value = 1.224035999999998
$ws.Cells.Item(14, 9).Value = 1.073352000000001
$ws.Cells.Item(14, 10).Value = 0.7055639999999994
$ws.Cells.Item(14, 11).Value = 0.892793999999999
$ws.Cells.Item(14, 12).Value = 0.9742109999999988
$ws.Cells.Item(14, 13).Value = 0.9902433333333326

# Row 15: Rotation-NoTilt
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 1.1
$ws.Cells.Item(15, 4).Value = 0.5
$ws.Cells.Item(15, 5).Value = 1.14
$ws.Cells.Item(15, 6).Value = 1.1
$ws.Cells.Item(15, 7).Value = 0.66
$ws.Cells.Item(15, 8).Value = 1.39
$ws.Cells.Item(15, 9).Value = 1.13
$ws.Cells.Item(15, 10).Value = 0.5
$ws.Cells.Item(15, 11).Value = 0.82
$ws.Cells.Item(15, 12).Value = 0.9600000000000001
$ws.Cells.Item(15, 13).Value = 0.9866666666666667

# Row 16: Rotation-60detTilt
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 1.05562511360001
$ws.Cells.Item(16, 4).Value = 0.7056244736000032
$ws.Cells.Item(16, 5).Value = 1.080000153599997
$ws.Cells.Item(16, 6).Value = 1.05562511360001
$ws.Cells.Item(16, 7).Value = 0.8029163007999992
$ws.Cells.Item(16, 8).Value = 1.223958758400002
$ws.Cells.Item(16, 9).Value = 1.073333478399991
$ws.Cells.Item(16, 10).Value = 0.7056244736000032
$ws.Cells.Item(16, 11).Value = 0.8928123136000001
$ws.Cells.Item(16, 12).Value = 0.9742187136000051
$ws.Cells.Item(16, 13).Value = 0.9902430464000004

# Row 17: HexGrid-90degTilt5degRes
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.9937378308183796
$ws.Cells.Item(17, 4).Value = 0.9960976807022145
$ws.Cells.Item(17, 5).Value = 0.9933726521081805
$ws.Cells.Item(17, 6).Value = 0.9937378308183796
$ws.Cells.Item(17, 7).Value = 0.9931397431958763
$ws.Cells.Item(17, 8).Value = 0.9943553639304488
$ws.Cells.Item(17, 9).Value = 0.9948838210213059
$ws.Cells.Item(17, 10).Value = 0.9960976807022145
$ws.Cells.Item(17, 11).Value = 0.9947351664051975
$ws.Cells.Item(17, 12).Value = 0.9942364986117886
$ws.Cells.Item(17, 13).Value = 0.9942645152960675

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 0.9895485387022636
$ws.Cells.Item(18, 4).Value = 1.014957106629773
$ws.Cells.Item(18, 5).Value = 0.9892628970588851
$ws.Cells.Item(18, 6).Value = 0.9895485387022636
$ws.Cells.Item(18, 7).Value = 1.002922244746067
$ws.Cells.Item(18, 8).Value = 0.9827887225805435
$ws.Cells.Item(18, 9).Value = 0.9897041957005024
$ws.Cells.Item(18, 10).Value = 1.014957106629773
$ws.Cells.Item(18, 11).Value = 1.002110001844329
$ws.Cells.Item(18, 12).Value = 0.9958292702732963
$ws.Cells.Item(18, 13).Value = 0.9948639509030057

# Row 19: HexGrid-60degTilt5degRes
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9829640018726238
$ws.Cells.Item(19, 4).Value = 1.052190682863402
$ws.Cells.Item(19, 5).Value = 0.9799104052442996
$ws.Cells.Item(19, 6).Value = 0.9829640018726238
$ws.Cells.Item(19, 7).Value = 1.027655813377856
$ws.Cells.Item(19, 8).Value = 0.9507899477795725
$ws.Cells.Item(19, 9).Value = 0.9792397958257774
$ws.Cells.Item(19, 10).Value = 1.052190682863402
$ws.Cells.Item(19, 11).Value = 1.016050544053851
$ws.Cells.Item(19, 12).Value = 0.9995072729632373
$ws.Cells.Item(19, 13).Value = 0.9954584411605886
